# Translate table header row (A1:H1) to Russian
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Название"
$ws.Range("B1").Value = "Квартира"
$ws.Range("C1").Value = "Расположение"
$ws.Range("D1").Value = "Жилой комплекс"
$ws.Range("E1").Value = "Адрес"
$ws.Range("F1").Value = "Цена"
$ws.Range("G1").Value = "Цена за м²"
$ws.Range("H1").Value = "Застройщик"
